$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Main")

$ws.Range("B5").Value = 7
$ws.Range("B6").Value = 3
$ws.Range("B7").Value = "J"
